# Weekly update: a new sampling date's data is inserted as the new row 11,
# pushing the existing data rows (11-134) down by one (to 12-135).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 11 (shifts 11..134 -> 12..135,
# and copies formatting from the row above, so the date style on column D
# carries through automatically).
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with this week's record.
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = 44490
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 100112006
$ws.Range("G11").Value = "Repollo"
$ws.Range("H11").Value = "Crespo record"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 180
$ws.Range("K11").Value = 600
$ws.Range("L11").Value = 700
$ws.Range("M11").Value = 650
$ws.Range("N11").Value = "$/unidad"
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 650
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"
